$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '26.615.54'
$ws.Range('E2').Value2 = '  +0.41%  '
$ws.Range('D3').Value2 = '1.738.57'
$ws.Range('E3').Value2 = '  +0.62%  '
$ws.Range('D4').Value2 = '''0.9997'
$ws.Range('E4').Value2 = '  +0.03%  '
$ws.Range('D5').Value2 = '''245.92'
$ws.Range('E5').Value2 = '  +0.40%  '
$ws.Range('D6').Value2 = '''1.001'
$ws.Range('E6').Value2 = '  +0.06%  '
$ws.Range('D7').Value2 = '''0.4967'
$ws.Range('E7').Value2 = '  +3.40%  '
$ws.Range('D8').Value2 = '''0.2674'
$ws.Range('E8').Value2 = '  +0.14%  '
$ws.Range('D9').Value2 = '''0.06265'
$ws.Range('E9').Value2 = '  +0.73%  '
$ws.Range('D10').Value2 = '1.746.83'
$ws.Range('E10').Value2 = '  +1.14%  '
$ws.Range('D11').Value2 = '''0.07031'
$ws.Range('E11').Value2 = '  -1.69%  '
$ws.Range('D12').Value2 = '''15.78'
$ws.Range('E12').Value2 = '  +0.70%  '
$ws.Range('D13').Value2 = '''4.600'
$ws.Range('E13').Value2 = '  +1.77%  '
$ws.Range('D14').Value2 = '''0.6128'
$ws.Range('E14').Value2 = '  -0.83%  '
$ws.Range('D15').Value2 = '''78.09'
$ws.Range('E15').Value2 = '  +1.19%  '
$ws.Range('D16').Value2 = '''1.001'
$ws.Range('E16').Value2 = '  +0.04%  '
$ws.Range('D17').Value2 = '26.616.77'
$ws.Range('E17').Value2 = '  +0.41%  '
$ws.Range('D18').Value2 = '''1.001'
$ws.Range('E18').Value2 = '  +0.17%  '
$ws.Range('D19').Value2 = '''0.000007250'
$ws.Range('E19').Value2 = '  +4.59%  '
$ws.Range('D20').Value2 = '''11.55'
$ws.Range('E20').Value2 = '  -0.93%  '
$ws.Range('D21').Value2 = '1.969.06'
$ws.Range('E21').Value2 = '  +1.03%  '
$ws.Range('D22').Value2 = '''4.546'
$ws.Range('E22').Value2 = '  +0.37%  '
$ws.Range('D23').Value2 = '''8.731'
$ws.Range('E23').Value2 = '  -2.50%  '
$ws.Range('D24').Value2 = '''5.291'
$ws.Range('E24').Value2 = '  +0.30%  '
$ws.Range('D25').Value2 = '''139.21'
$ws.Range('E25').Value2 = '  +2.26%  '
$ws.Range('D26').Value2 = '''15.42'
$ws.Range('E26').Value2 = '  +0.61%  '
$ws.Range('D27').Value2 = '''1.415'
$ws.Range('E27').Value2 = '  +0.58%  '
$ws.Range('B28').Value2 = 'BitcoinCash'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').Value2 = '''107.22'
$ws.Range('E28').Value2 = '  +0.51%  '
$ws.Range('B29').Value2 = 'LidoDAOToken'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value2 = '''1.752'
$ws.Range('E29').Value2 = '  -2.66%  '
$ws.Range('D30').Value2 = '''4.029'
$ws.Range('E30').Value2 = '  +1.25%  '
$ws.Range('D31').Value2 = '''0.08020'
$ws.Range('E31').Value2 = '  +0.00%  '
$ws.Range('D32').Value2 = '''3.731'
$ws.Range('E32').Value2 = '  +0.67%  '
$ws.Range('D33').Value2 = '''0.04596'
$ws.Range('E33').Value2 = '  +0.77%  '
$ws.Range('D34').Value2 = '''0.9998'
$ws.Range('E34').Value2 = '  +0.01%  '
$ws.Range('D35').Value2 = '''2.617'
$ws.Range('E35').Value2 = '  +0.10%  '
$ws.Range('D36').Value2 = '''1.015'
$ws.Range('E36').Value2 = '  +2.46%  '
$ws.Range('D37').Value2 = '''0.6365'
$ws.Range('E37').Value2 = '  +0.13%  '
$ws.Range('D38').Value2 = '''0.9053'
$ws.Range('E38').Value2 = '  -2.98%  '
$ws.Range('D39').Value2 = '''2.049'
$ws.Range('E39').Value2 = '  -2.27%  '
$ws.Range('E40').Value2 = '  +0.55%  '
$ws.Range('E41').Value2 = '  -0.32%  '
$ws.Range('D42').Value2 = '''0.01507'
$ws.Range('E42').Value2 = '  +0.26%  '
$ws.Range('D43').Value2 = '''101.68'
$ws.Range('E43').Value2 = '  -2.81%  '
$ws.Range('D44').Value2 = '''5.468'
$ws.Range('E44').Value2 = '  -3.27%  '
$ws.Range('D45').Value2 = '''0.3932'
$ws.Range('E45').Value2 = '  +0.46%  '
$ws.Range('D46').Value2 = '''6.866'
$ws.Range('E46').Value2 = '  -0.66%  '
$ws.Range('D47').Value2 = '''0.1176'
$ws.Range('E47').Value2 = '  -0.75%  '
$ws.Range('D48').Value2 = '''0.05386'
$ws.Range('E48').Value2 = '  +1.05%  '
$ws.Range('D49').Value2 = '''30.73'
$ws.Range('E49').Value2 = '  -0.77%  '
$ws.Range('D50').Value2 = '''7.816'
$ws.Range('E50').Value2 = '  -0.66%  '
$ws.Range('D51').Value2 = '''1.251'
$ws.Range('E51').Value2 = '  -1.09%  '
